$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that was bumped by one day
# (45243 -> 45244, i.e. 2023-11-13 -> 2023-11-14) for every data row (2-158).
$lastRow = 158
$range = $ws.Range("C2:C$lastRow")
$range.Value = 45244
